$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.503.10"
$ws.Range("E2").Value = "  +0.12%  "

$ws.Range("D3").Value = "1.952.57"
$ws.Range("E3").Value = "  +0.61%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.24"
$ws.Range("E5").Value = "  -0.13%  "

$ws.Range("E6").Value = "  +0.41%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.47"
$ws.Range("E7").Value = "  +5.64%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.377"
$ws.Range("E9").Value = "  +4.44%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0786"
$ws.Range("E10").Value = "  -7.82%  "

$ws.Range("E11").Value = "  +0.27%  "

$ws.Range("E12").Value = "  +6.33%  "

$ws.Range("D13").Value = "2.237.16"
$ws.Range("E13").Value = "  +0.48%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.825"
$ws.Range("E14").Value = "  +2.01%  "

$ws.Range("E15").Value = "  +1.18%  "

$ws.Range("E16").Value = "  +1.58%  "

$ws.Range("D17").Value = "1.944.67"
$ws.Range("E17").Value = "  +0.22%  "

$ws.Range("D18").Value = "36.434.82"
$ws.Range("E18").Value = "  -0.01%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.30"
$ws.Range("E19").Value = "  +0.00%  "

$ws.Range("D20").Value = "0.0₃0849"
$ws.Range("E20").Value = "  -1.61%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.07"
$ws.Range("E21").Value = "  +1.95%  "

$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.25"
$ws.Range("E22").Value = "  +0.44%  "

$ws.Range("E23").Value = "  +0.15%  "

$ws.Range("E24").Value = "  +3.74%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.36"
$ws.Range("E25").Value = "  +2.79%  "

$ws.Range("E26").Value = "  +7.73%  "

$ws.Range("E27").Value = "  -0.32%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "159.84"
$ws.Range("E28").Value = "  -0.59%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.27"
$ws.Range("E29").Value = "  +0.77%  "

$ws.Range("E30").Value = "  +20.78%  "

$ws.Range("E31").Value = "  +1.08%  "

$ws.Range("E32").Value = "  +3.41%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0612"
$ws.Range("E33").Value = "  -0.12%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.44"
$ws.Range("E34").Value = "  +6.30%  "

$ws.Range("E35").Value = "  +10.17%  "

$ws.Range("E36").Value = "  +0.14%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.27"
$ws.Range("E37").Value = "  +4.86%  "

$ws.Range("E38").Value = "  -1.45%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.45"
$ws.Range("E39").Value = "  -12.49%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0964"
$ws.Range("E40").Value = "  -1.69%  "

$ws.Range("E41").Value = "  +0.52%  "

$ws.Range("E42").Value = "  +2.06%  "

$ws.Range("E43").Value = "  +0.39%  "

$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.68"
$ws.Range("E44").Value = "  -1.12%  "

$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "1.358.41"
$ws.Range("E45").Value = "  +1.46%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.46"
$ws.Range("E46").Value = "  +2.63%  "

$ws.Range("E47").Value = "  -0.30%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.11"
$ws.Range("E48").Value = "  -1.01%  "

$ws.Range("E49").Value = "  +0.56%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "45.40"
$ws.Range("E50").Value = "  +5.09%  "

$ws.Range("D51").Value = "2.133.46"
$ws.Range("E51").Value = "  +0.71%  "
